# Trade #23 closed at 2026-02-16 22:54:22 - base_strategy UP +0.000%
# Append a new trade row (row 24) to both the "All Trades" sheet and the
# "base_strategy" sheet, mirroring the columns of the previous row (23).

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Duplicate the last row's structure (keeps the two "blank" cells -
    # Exit Price / Exit Reason - present as empty cells instead of being
    # dropped outright) into the new row, then overwrite with the new
    # trade's data.
    $ws.Range("A23:Q23").Copy($ws.Range("A24:Q24"))

    $ws.Cells.Item(24, 1).Value = 23

    $ws.Cells.Item(24, 2).NumberFormat = "@"
    $ws.Cells.Item(24, 2).Value = "2026-02-16"
    $ws.Cells.Item(24, 2).Style = "Normal"

    $ws.Cells.Item(24, 3).NumberFormat = "@"
    $ws.Cells.Item(24, 3).Value = "22:54:22"
    $ws.Cells.Item(24, 3).Style = "Normal"

    $ws.Cells.Item(24, 4).Value = "base_strategy"
    $ws.Cells.Item(24, 5).Value = "UP"
    $ws.Cells.Item(24, 6).Value = 49.999998
    # Column G (Exit Price) stays blank - left as copied.
    $ws.Cells.Item(24, 8).Value = "OPEN"
    $ws.Cells.Item(24, 9).Value = 0
    $ws.Cells.Item(24, 10).Value = 0
    $ws.Cells.Item(24, 11).Value = 100
    $ws.Cells.Item(24, 12).Value = 0
    $ws.Cells.Item(24, 13).Value = 0
    $ws.Cells.Item(24, 14).Value = 0.6
    $ws.Cells.Item(24, 15).Value = "Normal spread capture: 19600 bps"
    # Column P (Exit Reason) stays blank - left as copied.
    $ws.Cells.Item(24, 17).Value = 0
}
